$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: scroll / "Ships by" gains locatorType "p" and waitBefore count "1" ---
$ws.Range("E3").Value = "p"
$ws.Range("F3").Value = 1

# --- Row 4: elementClick -> click ---
$ws.Range("C4").Value = "click"

# --- Row 6: wairfortext -> waitfortext, gains locatorType "h2" + wait timings ---
$ws.Range("C6").Value = "waitfortext"
$ws.Range("E6").Value = "h2"
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = 2000

# --- Row 7: scroll / "We Accept" gains locatorType "p" and waitBefore count "1" ---
$ws.Range("E7").Value = "p"
$ws.Range("F7").Value = 1

# --- Row 8: elementClick -> click, waitAfter 5000 -> 9000 ---
$ws.Range("C8").Value = "click"
$ws.Range("H8").Value = 9000

# --- Row 9: wairfortext -> waitfortext, text changes, gains locatorType "h2" + wait timings ---
$ws.Range("C9").Value = "waitfortext"
$ws.Range("D9").Value = "Contact information"
$ws.Range("E9").Value = "h2"
$ws.Range("G9").Value = 1000
$ws.Range("H9").Value = 2000

# --- Row 10: elementClick -> click ---
$ws.Range("C10").Value = "click"

# --- Row 11: Enabled Yes -> no, elementinputdata -> type ---
$ws.Range("B11").Value = "no"
$ws.Range("C11").Value = "type"

# --- Row 12: Enabled Yes -> no, elementClick -> click ---
$ws.Range("B12").Value = "no"
$ws.Range("C12").Value = "click"

# --- Row 13: Enabled Yes -> no, wairfortext -> waitfortext, text changes, gains locatorType "h3" + wait timings ---
$ws.Range("B13").Value = "no"
$ws.Range("C13").Value = "waitfortext"
$ws.Range("D13").Value = "Email Address"
$ws.Range("E13").Value = "h3"
$ws.Range("G13").Value = 1000
$ws.Range("H13").Value = 5000

# --- Row 14: Enabled Yes -> no ---
$ws.Range("B14").Value = "no"

# --- New row 16: formatted-but-empty placeholder cells (light grey 8pt Segoe UI) ---
$r16 = $ws.Range("E16:F16")
$r16.Style = "Normal"
$r16.Font.Name = "Segoe UI"
$r16.Font.Size = 8
$r16.Font.Color = 13421772

# --- Update the visible selection to match the author's last selection ---
$ws.Range("B11:B14").Select()
